$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2 changes from 3.75 to 3.7; the old D2/E2/F2 values are relocated to I2/J2/K2,
# and E2/F2 get new values.
$ws.Range("D2").Value = 3.7
$ws.Range("E2").Value = "12.97, 18.8"
$ws.Range("F2").Formula = "'0.25"
$ws.Range("I2").Value = 3.75
$ws.Range("J2").Value = "13.02, 18.75"
$ws.Range("K2").Formula = "'0.3"

# Row 3: E3 keeps the same text but via its (re-indexed) shared string.
$ws.Range("E3").Value = "6.46, 9.51"

# Row 1: new label cell.
$ws.Range("J1").Value = "used to be"

# Update selection to match the target workbook state.
$ws.Range("E2").Select() | Out-Null
